$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 65, pushing existing rows 65-129 down to 66-130.
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new data record.
$ws.Range("A65").Value = 7
$ws.Range("B65").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C65").Value = "Ñuble"
$ws.Range("D65").Value = 45175
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E65").Value = 16
$ws.Range("F65").Value = 100112013
$ws.Range("G65").Value = "Alcachofa"
$ws.Range("H65").Value = "Madrigal"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 80
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = 13000
$ws.Range("N65").Value = "`$/caja 40 unidades"
$ws.Range("O65").Value = "Provincia de Limarí"
$ws.Range("P65").Value = 325
$ws.Range("Q65").Value = 40
$ws.Range("R65").Value = "Hortaliza"
